$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 18: full data row for the 20230831 session ---
$ws.Cells.Item(18, 1).Value = 20230831
$ws.Cells.Item(18, 2).Value = 2
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 4).Value = 4
$ws.Cells.Item(18, 5).Value = 5
$ws.Cells.Item(18, 6).Value = 6
$ws.Cells.Item(18, 7).Value = 7

# --- New row 19: thresholded/error values for that session ---
$ws.Cells.Item(19, 2).Value = 18
$ws.Cells.Item(19, 3).Value = "8,14,15,17,20"
$ws.Cells.Item(19, 4).Value = "10,17,18"
$ws.Cells.Item(19, 5).Value = " 4,8,13"
$ws.Cells.Item(19, 6).Value = "17,18"

# --- Mark the A11 (20230817) session cell the same way A9 already is (quote-prefixed style) ---
$ws.Range("A11").Style = $ws.Range("A9").Style

# --- Make that shared font style red (applies to both A9 and A11 now) ---
$ws.Range("A9").Font.Color = 255

# --- Update the selection/active cell to reflect where the user ended up ---
$ws.Application.Goto($ws.Range("C20"))

$wb.Save()
